# [ADDITIONAL SCRAPING] added scraping code for extra bowling attributes and excel sheets
#
# 1. Adds a new worksheet "ODI Bowling Extra" (after "ODI Batting Extra") holding
#    MATCH_CODE / MAIDEN_OVERS / PERCENT_WICKETS_OF_ALL scraped stats.
# 2. Tidies up the existing "ODI Batting Extra" sheet by dropping the leftover
#    empty placeholder cells that the old scraper used to emit.

function Set-TextCell($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    # Force text storage so numeric-looking strings ("10", "10.00%", match
    # codes, ...) are not silently reinterpreted as numbers/percentages.
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. New sheet: "ODI Bowling Extra"
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "ODI Bowling Extra"

# Header row (bold, bordered, centered/top-aligned - matches the other
# "Extra" sheets' header styling)
$header = $newSheet.Range("A1:C1")
$header.Font.Bold = $true
$header.Borders.LineStyle = 1
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160

$newSheet.Range("A1").Value = "MATCH_CODE"
$newSheet.Range("B1").Value = "MAIDEN_OVERS"
$newSheet.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"

Set-TextCell $newSheet "A2" "4172"
Set-TextCell $newSheet "B2" "1"
Set-TextCell $newSheet "C2" "10.00%"
Set-TextCell $newSheet "A3" "4174"
Set-TextCell $newSheet "B3" "0"
Set-TextCell $newSheet "A4" "4178"
Set-TextCell $newSheet "A5" "4194"
Set-TextCell $newSheet "B5" "1"
Set-TextCell $newSheet "A6" "4197"
Set-TextCell $newSheet "B6" "1"
Set-TextCell $newSheet "A7" "4201"
Set-TextCell $newSheet "B7" "0"
Set-TextCell $newSheet "A8" "4241"
Set-TextCell $newSheet "B8" "0"
Set-TextCell $newSheet "A9" "4244"
Set-TextCell $newSheet "B9" "1"
Set-TextCell $newSheet "C9" "10.00%"
Set-TextCell $newSheet "A10" "4247"
Set-TextCell $newSheet "A11" "4273"
Set-TextCell $newSheet "B11" "0"
Set-TextCell $newSheet "A12" "4304"
Set-TextCell $newSheet "B12" "0"
Set-TextCell $newSheet "C12" "30.00%"
Set-TextCell $newSheet "A13" "4308"
Set-TextCell $newSheet "A14" "4319"
Set-TextCell $newSheet "B14" "2"
Set-TextCell $newSheet "C14" "50.00%"
Set-TextCell $newSheet "A15" "4324"
Set-TextCell $newSheet "B15" "1"
Set-TextCell $newSheet "C15" "30.00%"
Set-TextCell $newSheet "A16" "4334"
Set-TextCell $newSheet "B16" "1"
Set-TextCell $newSheet "C16" "20.00%"
Set-TextCell $newSheet "A17" "4337"
Set-TextCell $newSheet "A18" "4340"
Set-TextCell $newSheet "B18" "1"
Set-TextCell $newSheet "A19" "4349"
Set-TextCell $newSheet "B19" "0"
Set-TextCell $newSheet "C19" "10.00%"
Set-TextCell $newSheet "A20" "4375"
Set-TextCell $newSheet "B20" "1"
Set-TextCell $newSheet "C20" "10.00%"
Set-TextCell $newSheet "A21" "4376"
Set-TextCell $newSheet "B21" "0"
Set-TextCell $newSheet "C21" "30.00%"

# ---------------------------------------------------------------------------
# 2. Clean up "ODI Batting Extra": clear the old empty placeholder cells
#    (rows where NUM_4 / NUM_6 / PERCENT_RUNS_OF_TOTAL / ... had no data but
#    still carried an empty inline-string cell).
# ---------------------------------------------------------------------------
$battingExtra = $wb.Worksheets.Item("ODI Batting Extra")
$emptyRanges = $battingExtra.Range("C2:E2,B3:E3,C4:E4,C6:E6,C7:E7,C8:E8,B9:E9,C10:E10,B11:E11,B13:E13,E14,C15:E15,B16:E16,C17:E17,C19:E19,C20:E20,B21:F21")
foreach ($area in $emptyRanges.Areas) {
    $area.Value = ""
}
